{"js": "// COV_motion.docx \u2014 certificate-of-service \"else\" branch tweak.\n//\n// The template's fallback (non-cert_of_service) sentence currently reads:\n//   \"else %}______________, ______________________________\"\n// It must match the \"if\" branch's phrasing (which already has \", I, \")\n// by inserting \" I,\" right after the first blank/comma, i.e. becoming:\n//   \"else %}______________, I, ______________________________\"\n\nconst body = context.document.body;\n\nconst oldText = \"else %}______________, ______________________________\";\nconst results = body.search(oldText, { matchCase: true, matchWildcards: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the target 'else %}' cert-of-service fallback text.\");\n}\n\n// Within the matched range, locate the comma that follows the first blank\n// line, then insert \" I,\" immediately after it (keeping the rest of the\n// sentence, including the leading space before the long blank, intact).\nconst target = results.items[0];\nconst commaMatches = target.search(\",\", { matchCase: true, matchWildcards: false });\ncommaMatches.load(\"items\");\nawait context.sync();\n\nif (commaMatches.items.length === 0) {\n  throw new Error(\"Could not find the comma to insert ' I,' after.\");\n}\n\ncommaMatches.items[0].insertText(\" I,\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# COV_motion.docx \u2014 certificate-of-service \"else\" branch tweak.\n#\n# The template's fallback (non-cert_of_service) sentence currently reads:\n#   \"else %}______________, ______________________________\"\n# It must match the \"if\" branch's phrasing (which already has \", I, \")\n# by inserting \" I,\" right after the first blank/comma, i.e. becoming:\n#   \"else %}______________, I, ______________________________\"\n\n$d = $word.ActiveDocument\n\n$oldText = \"else %}______________, ______________________________\"\n$newText = \"else %}______________, I, ______________________________\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.Forward = $true\n$find.Wrap = 1\n\n# wdFindContinue=1, wdReplaceAll=2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
